$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert the new first paragraph (title + italic word-problem body)
#    at the very start of the document. InsertXML on a collapsed range
#    splices in whole new paragraph(s) before the existing content,
#    leaving the pre-existing paragraph (and its _GoBack bookmark)
#    completely untouched right after it.
# ---------------------------------------------------------------------
$insertPoint = $d.Range(0, 0)

$newParaXml = '<?xml version="1.0" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body>' + `
'<w:p>' + `
'<w:pPr><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Times New Roman"/></w:rPr></w:pPr>' + `
'<w:r><w:rPr><w:b/><w:color w:val="365F91" w:themeColor="accent1" w:themeShade="BF"/></w:rPr><w:t>A Cat, a Parrot, and a Bag of Seed:</w:t></w:r>' + `
'<w:r><w:rPr><w:b/><w:color w:val="365F91" w:themeColor="accent1" w:themeShade="BF"/></w:rPr><w:br/></w:r>' + `
'<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:i/></w:rPr><w:t>A man finds himself on a riverbank with a cat, a parrot and a bag of seed. He needs to transport all three to the other side of the river in his boat. However, the boat has room for only the man himself and one other item (either the cat, parrot or seed). In his absence, the cat could eat the parrot, and the parrot would eat the bag of seed. Show how he can get all the passengers to the other side, without leaving the wrong ones alone together.</w:t></w:r>' + `
'</w:p>' + `
'</w:body></w:document>' + `
'</pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($newParaXml)

# ---------------------------------------------------------------------
# 2) The original (now second) paragraph still holds the _GoBack
#    bookmark. Replace its full range (including its own paragraph
#    mark) with paragraph-mark formatting (bold + themed color) and a
#    trailing bold line-break run. Because the bookmark has zero width
#    at the start of that range, it is preserved as-is (same id/name)
#    rather than duplicated, as long as we do not redeclare it in the
#    replacement fragment.
# ---------------------------------------------------------------------
$secondPara = $d.Paragraphs.Item(2)
$fullSecond = $d.Range($secondPara.Range.Start, $secondPara.Range.End)

$secondParaXml = '<?xml version="1.0" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body>' + `
'<w:p>' + `
'<w:pPr><w:rPr><w:b/><w:color w:val="365F91" w:themeColor="accent1" w:themeShade="BF"/></w:rPr></w:pPr>' + `
'<w:r><w:rPr><w:b/></w:rPr><w:br/></w:r>' + `
'</w:p>' + `
'</w:body></w:document>' + `
'</pkg:xmlData></pkg:part></pkg:package>'

$fullSecond.InsertXML($secondParaXml)
